$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: fiber-sedan -> turn into a "Proposal" line item with price reset to 0/empty
$ws.Range("B4").Value = "Fiber Optic Lights Sedan Proposal "
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = ""

# Row 6: fiber-suv -> turn into a "Proposal" line item with price reset to 0/empty
$ws.Range("B6").Value = "Fiber Optic Lights SUV Proposal "
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = ""

# Update the current selection to match the saved view (B6 selected, scrolled back to top)
$ws.Range("B6").Select()
